$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: 2/6/2020 guest-lecture entry -------------------------------
$ws.Range("A26").Value = 43867
$ws.Range("A26").NumberFormat = "m/d/yy"
$ws.Range("B26").Value = "5:00 - 7:50 pm"
$ws.Range("C26").Value = "Myself"
$ws.Range("D26").Value = "Learn about KEP4,5,6. Learn about mental simulations"
$ws.Range("E26").Value = "Understood how experts approach reading and analyzing code. Reflecting back is one of the most important things to do."
$ws.Range("F26").Value = "Reflecting back, applying key expert practices will come slowly with time. Templates can help us achieve that slowly and steadily"
$ws.Range("G26").Value = "Excited to listen to the guest lecture, and learn about these key practices."
$ws.Rows(26).RowHeight = 68

# --- Row 27: 2/10/2020 midterm-study entry ------------------------------
$ws.Range("A27").Value = 43871
$ws.Range("A27").NumberFormat = "m/d/yy"
$ws.Range("B27").Value = "9:00 pm - 10:00pm"
$ws.Range("C27").Value = "Myself"
$ws.Range("D27").Value = "Study a little bit for the mid term"
$ws.Range("E27").Value = "Covered slides 1-3"
$ws.Range("F27").Value = "The concepts are easy to read, but I understand their application will come over time"
$ws.Range("G27").Value = "Neutral"
$ws.Rows(27).RowHeight = 51

# --- Row 28: 2/11 & 2/12/2020 more midterm-study entry ------------------
$ws.Range("B28").Value = "9:00pm - 11:00pm"
$ws.Range("C28").Value = "Myself"
$ws.Range("D28").Value = "Study more for the mid term"
$ws.Range("E28").Value = "Covered slides 4 and 5 and UML notations"
$ws.Range("A28").Value = 43867
$ws.Range("A28").NumberFormat = "m/d/yy"
$ws.Range("A28").Value = "2/11/2020 and 2/12/2020"
$ws.Range("F28").Value = "Concepts in the latter slides are a little less easy to grasp, but a revision would reinforce these"
$ws.Range("G28").Value = "Neutral, a little tensed"
$ws.Rows(28).RowHeight = 51

# --- Selection matches the committed view state -------------------------
$ws.Range("G28").Select()
